# Daman's paragraph stays as-is, but the "_GoBack" bookmark that used to
# sit right after it needs to move down onto the new paragraph Manoj is
# adding. Remove the old (now-misplaced) bookmark first so there's no
# name clash once the new one is created below.
$d = $word.ActiveDocument

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# Append the new paragraph at the very end of the document. We build it
# from raw WordprocessingML so we can get the exact run/proofErr layout
# Word itself would produce for a flagged-as-misspelled word ("Manoj"
# wrapped in spellStart/spellEnd, split into its own run) and so the
# _GoBack bookmark lands collapsed right after the last run (the usual
# "last edit position" spot), instead of wrapping the whole paragraph.
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

$newParagraphXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Manoj</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> edited second thing.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newParagraphXml)
